# Apply weekly update: insert a new data row at row 42 (pushing existing
# rows 42..99 down to 43..100) and populate the new row with this week's
# Alcachofa price record for "Macroferia Regional de Talca".

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at position 42; all rows from 42 downward shift
# down by one (old row 42 -> new row 43, ..., old row 99 -> new row 100).
$ws.Rows(42).Insert()

# Populate the newly inserted row 42 with the new record's values.
$ws.Cells.Item(42, 1).Value = 5
$ws.Cells.Item(42, 2).Value = "Macroferia Regional de Talca"
$ws.Cells.Item(42, 3).Value = "Maule"
$ws.Cells.Item(42, 4).Value = 44799
$ws.Cells.Item(42, 5).Value = 7
$ws.Cells.Item(42, 6).Value = 100112013
$ws.Cells.Item(42, 7).Value = "Alcachofa"
$ws.Cells.Item(42, 8).Value = "Madrigal"
$ws.Cells.Item(42, 9).Value = "Primera"
$ws.Cells.Item(42, 10).Value = 300
$ws.Cells.Item(42, 11).Value = 12000
$ws.Cells.Item(42, 12).Value = 12000
$ws.Cells.Item(42, 13).Value = 12000
$ws.Cells.Item(42, 14).Value = "`$/caja 40 unidades"
$ws.Cells.Item(42, 15).Value = "Provincia del Elquí"
$ws.Cells.Item(42, 16).Value = 300
$ws.Cells.Item(42, 17).Value = 40
$ws.Cells.Item(42, 18).Value = "Hortaliza"
